# Git batch file commit June 03
#
# - Renames the shared "selenium3888223" username value (used on both the
#   Register sheet's "newUserName" column and the LogIn sheet's userName
#   column) to "newUser34422".
# - Switches the active tab from "Register" to "LogIn" and updates each
#   sheet's remembered selection.

$wb = $excel.ActiveWorkbook

$wsRegister = $wb.Worksheets.Item("Register")
$wsLogin    = $wb.Worksheets.Item("LogIn")

# Update the shared user name value on both sheets that reference it so the
# underlying shared string is edited in place (rather than orphaning the old
# string and adding a new one).
$wsRegister.Range("I2").Value = "newUser34422"
$wsLogin.Range("A2").Value = "newUser34422"

# Register sheet: clear its old selection/scroll position.
$wsRegister.Range("I2").Select() | Out-Null

# Make LogIn the active sheet (this also flips workbookView's activeTab and
# moves tabSelected onto the LogIn sheet's view) and set its new selection.
$wsLogin.Activate()
$wsLogin.Range("B12").Select() | Out-Null
